$wb = $excel.ActiveWorkbook

# Sheets 1-4: Potencia Acumulada, Geracao Periodo Medio, Atendimento a Ponta,
# Potencia Incremental -> fix E1 label, remove Total row (row 13)
$sheetsToFix = @(
    @{ Index = 1; Label = "2050" },
    @{ Index = 2; Label = "2050" },
    @{ Index = 3; Label = "2050" },
    @{ Index = 4; Label = "2041-2050" }
)

foreach ($item in $sheetsToFix) {
    $ws = $wb.Worksheets.Item($item.Index)
    $ws.Range("E1").Value = $item.Label
    $ws.Rows.Item(13).Delete()
}

# Sheet 5: Emissoes Totais (MtCO2eq) -> only fix E1 label, no Total row present
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").Value = "2050"

# Sheet 6: Custo Total (bilhões de R$) -> remove Total row (row 4)
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
